$wb = $excel.ActiveWorkbook
$wsTest = $wb.Worksheets.Item(1)   # "Test Data"
$wsMeas = $wb.Worksheets.Item(2)   # "Measurement Data"

# --- "Test Data" sheet: new test run (Id 56 -> 61) with a fresh CreatedDate ---
$wsTest.Range("A2").Value = 61
$wsTest.Range("G2").Value = 44042.32261403935

# --- "Measurement Data" sheet: rows now belong to TestId 61, Ids shifted +100, ---
# --- refreshed CreatedDate timestamps, and re-measured Current reading deltas ---

# Row 2
$wsMeas.Range("A2").Value = 446
$wsMeas.Range("T2").Value = 3.38546
$wsMeas.Range("Y2").Value = 44042.32274178241
$wsMeas.Range("Z2").Value = 61

# Row 3
$wsMeas.Range("A3").Value = 447
$wsMeas.Range("T3").Value = 1.25649
$wsMeas.Range("Y3").Value = 44042.32285135417
$wsMeas.Range("Z3").Value = 61

# Row 4
$wsMeas.Range("A4").Value = 448
$wsMeas.Range("T4").Value = 0.05478
$wsMeas.Range("Y4").Value = 44042.32295204861
$wsMeas.Range("Z4").Value = 61

# Row 5
$wsMeas.Range("A5").Value = 449
$wsMeas.Range("T5").Value = 0.01692
$wsMeas.Range("Y5").Value = 44042.323061689814
$wsMeas.Range("Z5").Value = 61

# Row 6
$wsMeas.Range("A6").Value = 450
$wsMeas.Range("Y6").Value = 44042.32316165509
$wsMeas.Range("Z6").Value = 61

# Row 7
$wsMeas.Range("A7").Value = 451
$wsMeas.Range("U7").Value = 2.98076
$wsMeas.Range("Y7").Value = 44042.323287928244
$wsMeas.Range("Z7").Value = 61

# Row 8
$wsMeas.Range("A8").Value = 452
$wsMeas.Range("U8").Value = 1.49365
$wsMeas.Range("Y8").Value = 44042.32339849537
$wsMeas.Range("Z8").Value = 61

# Row 9
$wsMeas.Range("A9").Value = 453
$wsMeas.Range("U9").Value = 0.29824
$wsMeas.Range("Y9").Value = 44042.32349927083
$wsMeas.Range("Z9").Value = 61

# Row 10
$wsMeas.Range("A10").Value = 454
$wsMeas.Range("U10").Value = 0.07969
$wsMeas.Range("Y10").Value = 44042.32360899306
$wsMeas.Range("Z10").Value = 61

# Row 11
$wsMeas.Range("A11").Value = 455
$wsMeas.Range("U11").Value = 0.00086
$wsMeas.Range("Y11").Value = 44042.323708993055
$wsMeas.Range("Z11").Value = 61

# Row 12
$wsMeas.Range("A12").Value = 456
$wsMeas.Range("V12").Value = 2.9905
$wsMeas.Range("Y12").Value = 44042.32382549768
$wsMeas.Range("Z12").Value = 61

# Row 13
$wsMeas.Range("A13").Value = 457
$wsMeas.Range("V13").Value = 1.17426
$wsMeas.Range("Y13").Value = 44042.32393503472
$wsMeas.Range("Z13").Value = 61

# Row 14
$wsMeas.Range("A14").Value = 458
$wsMeas.Range("V14").Value = 0.09117
$wsMeas.Range("Y14").Value = 44042.32403521991
$wsMeas.Range("Z14").Value = 61

# Row 15
$wsMeas.Range("A15").Value = 459
$wsMeas.Range("V15").Value = 0.02457
$wsMeas.Range("Y15").Value = 44042.32413587963
$wsMeas.Range("Z15").Value = 61

# Row 16
$wsMeas.Range("A16").Value = 460
$wsMeas.Range("V16").Value = 0.00045
$wsMeas.Range("Y16").Value = 44042.324246412034
$wsMeas.Range("Z16").Value = 61

# Row 17
$wsMeas.Range("A17").Value = 461
$wsMeas.Range("W17").Value = 2.86945
$wsMeas.Range("Y17").Value = 44042.324371759256
$wsMeas.Range("Z17").Value = 61

# Row 18
$wsMeas.Range("A18").Value = 462
$wsMeas.Range("W18").Value = 1.23445
$wsMeas.Range("Y18").Value = 44042.32448113426
$wsMeas.Range("Z18").Value = 61

# Row 19
$wsMeas.Range("A19").Value = 463
$wsMeas.Range("W19").Value = 0.10529
$wsMeas.Range("Y19").Value = 44042.324581863424
$wsMeas.Range("Z19").Value = 61

# Row 20
$wsMeas.Range("A20").Value = 464
$wsMeas.Range("W20").Value = 0.02899
$wsMeas.Range("Y20").Value = 44042.32468194445
$wsMeas.Range("Z20").Value = 61

# Row 21
$wsMeas.Range("A21").Value = 465
$wsMeas.Range("W21").Value = 0.0004
$wsMeas.Range("Y21").Value = 44042.32479293981
$wsMeas.Range("Z21").Value = 61
